# "load list student UI" - add header rows to all data sheets so the
# Java UI can read column names (Name/Mssv/Sex/Cmnd/ClassId, Class Id/Class Name,
# Subject Id/Subject Name, Class Id/Subject ID/Room), plus minor view tweaks.

$wb = $excel.ActiveWorkbook

$wsStudents = $wb.Worksheets.Item(1)   # "students"
$wsClass    = $wb.Worksheets.Item(2)   # "class"
$wsSubject  = $wb.Worksheets.Item(3)   # "subject"
$wsSchedule = $wb.Worksheets.Item(4)   # "schedule"

# --- students sheet: insert header row, fill labels, widen a couple of columns
$wsStudents.Rows.Item(1).Insert()
$wsStudents.Range("A1").Value = "Name"
$wsStudents.Range("B1").Value = "Mssv"
$wsStudents.Range("C1").Value = "Sex"
$wsStudents.Range("D1").Value = "Cmnd"
$wsStudents.Range("E1").Value = "ClassId"
$wsStudents.Columns.Item(1).ColumnWidth = 24.0
$wsStudents.Columns.Item(4).ColumnWidth = 24.166666666666668

# --- class sheet: insert header row, fill labels
$wsClass.Rows.Item(1).Insert()
$wsClass.Range("A1").Value = "Class Id"
$wsClass.Range("B1").Value = "Class Name"

# --- subject sheet: insert header row, fill labels
$wsSubject.Rows.Item(1).Insert()
$wsSubject.Range("A1").Value = "Subject Id"
$wsSubject.Range("B1").Value = "Subject Name"

# --- schedule sheet: insert header row, fill labels
$wsSchedule.Rows.Item(1).Insert()
$wsSchedule.Range("A1").Value = "Class Id"
$wsSchedule.Range("B1").Value = "Subject ID"
$wsSchedule.Range("C1").Value = "Room"

# --- restore per-sheet selections
$wsStudents.Range("I20").Select() | Out-Null
$wsSubject.Range("H13").Select() | Out-Null
$wsSchedule.Range("C1").Select() | Out-Null

# "class" ends up the active/visible tab with its own selection
$wsClass.Activate() | Out-Null
$wsClass.Range("I14").Select() | Out-Null
